# Removed Type from imports
# The "Type" column (G) - together with its sample values ("Pool") - is no
# longer needed, so delete the whole column. Excel shifts every column to
# its right (Folio No, Instrument, Sector, Investment Domicile *,
# Custom Field 1) one position to the left, the shared-string table is
# compacted accordingly, and the selection lands on the column that now
# occupies the freed-up G slot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G").Delete()
$ws.Columns("G").Select()
